# Add a new "2022-Q1" sheet (fund-holding detail) before the "总计" summary
# sheet, and prepend a corresponding "2022-Q1" row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $value) {
    # Force a numeric-looking string (e.g. "17.74" or "002387") to be stored
    # as text rather than being auto-coerced into a Number by Excel, while
    # avoiding the left-over quote-prefix style that NumberFormat="@" alone
    # would leave behind.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet right after "2021-Q4" (i.e. right
#    before "总计"), cloning the column formatting of "2021-Q4".
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsNew = $wb.Worksheets.Add($null, $wsQ4)
$wsNew.Name = "2022-Q1"

# Copy just the formatting (no values) from "2021-Q4" - done in two pieces
# because row 1 has no A1 cell (column A only starts at row 2) in the
# source sheet, and a blanket A1:H12 copy would otherwise manufacture an
# empty <c r="A1"/> cell that doesn't exist in the reference layout.
$wsQ4.Range("B1:H12").Copy()
$wsNew.Range("B1").PasteSpecial(-4122)
$wsQ4.Range("A2:A12").Copy()
$wsNew.Range("A2").PasteSpecial(-4122)

# Header row
$wsNew.Cells.Item(1,2).Value = "基金代码"
$wsNew.Cells.Item(1,3).Value = "基金名称"
$wsNew.Cells.Item(1,4).Value = "基金规模"
$wsNew.Cells.Item(1,5).Value = "股票总仓位"
$wsNew.Cells.Item(1,6).Value = "仓位占比"
$wsNew.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsNew.Cells.Item(1,8).Value = "仓位排名"

$data = @(
    @("002387","工银瑞信沪港深股票A","17.74","94.04","4.23","0.7504",9),
    @("005644","广发沪港深行业龙头混合","13.85","88.27","4.84","0.6703",6),
    @("005583","易方达港股通红利灵活配置混合","6.88","87.72","5.57","0.3832",3),
    @("002653","泰康沪港深精选灵活配置混合","7.89","87.05","2.09","0.1649",9),
    @("003580","泰康沪港深价值优选灵活配置混合","1.80","85.70","2.05","0.0369",10),
    @("007512","工银瑞信沪港深股票C","0.86","94.04","4.23","0.0364",9),
    @("007132","长城港股通价值精选多策略混合","0.96","85.80","2.55","0.0245",9),
    @("000927","博时大中华亚太精选股票(QDII) - 美元现汇","0.32","92.94","3.89","0.0124",10),
    @("050015","博时大中华亚太精选股票(QDII) -人民币","0.32","92.94","3.89","0.0124",10),
    @("161124","易方达香港恒生综合小型股指数（QDII-LOF）A","0.28","92.62","1.79","0.0050",4),
    @("006263","易方达香港恒生综合小型股指数（QDII-LOF）C","0.06","92.62","1.79","0.0011",4)
)

$r = 2
foreach ($row in $data) {
    $wsNew.Cells.Item($r,1).Value = $r - 2
    Set-TextValue $wsNew.Cells.Item($r,2) $row[0]
    $wsNew.Cells.Item($r,3).Value = $row[1]
    Set-TextValue $wsNew.Cells.Item($r,4) $row[2]
    Set-TextValue $wsNew.Cells.Item($r,5) $row[3]
    Set-TextValue $wsNew.Cells.Item($r,6) $row[4]
    Set-TextValue $wsNew.Cells.Item($r,7) $row[5]
    $wsNew.Cells.Item($r,8).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing rows down by one and renumbering column A.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q1"
$wsTotal.Cells.Item(2,3).Value = 11
$wsTotal.Cells.Item(2,4).Value = 2.1

# Renumber the now-shifted rows (were 0,1,2,3 -> now 1,2,3,4)
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(6,1).Value = 4
